$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176 (shifts existing rows 176:294 down to 177:295)
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(176, 1).Value  = 10
$ws.Cells.Item(176, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(176, 3).Value  = "La Araucanía"
$ws.Cells.Item(176, 4).Value  = 44596
$ws.Cells.Item(176, 5).Value  = 9
$ws.Cells.Item(176, 6).Value  = "Fruta"
$ws.Cells.Item(176, 7).Value  = 100108
$ws.Cells.Item(176, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(176, 9).Value  = 100108002
$ws.Cells.Item(176, 10).Value = "Mango"
$ws.Cells.Item(176, 11).Value = "Sin especificar"
$ws.Cells.Item(176, 12).Value = "Primera"
$ws.Cells.Item(176, 13).Value = 315
$ws.Cells.Item(176, 14).Value = 7500
$ws.Cells.Item(176, 15).Value = 8000
$ws.Cells.Item(176, 16).Value = 7754
$ws.Cells.Item(176, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(176, 18).Value = "Perú"
$ws.Cells.Item(176, 19).Value = 1938
$ws.Cells.Item(176, 20).Value = 4
